# Weekly update: Hortaliza, Vega Monumental Concepción - Melón
#
# Updates the existing week (rows 173-178, Calameño/Tuna x Extra/Primera/Segunda,
# "Región de O'Higgins") with a new date & revised prices, and re-adds the
# previous "Región de O'Higgins" week (old rows 179-184) plus the always-kept
# "Calameño Primera" outlier record (old row 179) and the "Región Metropolitana"
# week (old rows 180-185) as new rows, extending the sheet from 185 to 191 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmt = "YYYY-MM-DD HH:MM:SS"

function Set-DataRow {
    param(
        [int]$Row,
        [double]$Date,
        [string]$Variedad,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [string]$Origen,
        [double]$PrecioKg
    )

    $ws.Cells.Item($Row, 1).Value = 11
    $ws.Cells.Item($Row, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($Row, 3).Value = "Bíobío"
    $ws.Cells.Item($Row, 4).Value = $Date
    $ws.Cells.Item($Row, 4).NumberFormat = $dateFmt
    $ws.Cells.Item($Row, 5).Value = 8
    $ws.Cells.Item($Row, 6).Value = 100112027
    $ws.Cells.Item($Row, 7).Value = "Melón"
    $ws.Cells.Item($Row, 8).Value = $Variedad
    $ws.Cells.Item($Row, 9).Value = $Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = "$/unidad"
    $ws.Cells.Item($Row, 15).Value = $Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = 1
    $ws.Cells.Item($Row, 18).Value = "Hortaliza"
}

# Rows 173-178: same combos as before, date bumped 44194 -> 44568, prices raised.
Set-DataRow 173 44568 "Calameño" "Extra" 1000 900 900 900 "Región de O'Higgins" 900
Set-DataRow 174 44568 "Calameño" "Primera" 1000 800 800 800 "Región de O'Higgins" 800
Set-DataRow 175 44568 "Calameño" "Segunda" 1000 700 700 700 "Región de O'Higgins" 700
Set-DataRow 176 44568 "Tuna" "Extra" 1000 900 900 900 "Región de O'Higgins" 900
Set-DataRow 177 44568 "Tuna" "Primera" 1000 800 800 800 "Región de O'Higgins" 800
Set-DataRow 178 44568 "Tuna" "Segunda" 1000 700 700 700 "Región de O'Higgins" 700

# Rows 179-184: re-introduce the older "Región de O'Higgins" week (date 44194).
Set-DataRow 179 44194 "Calameño" "Extra" 400 800 800 800 "Región de O'Higgins" 800
Set-DataRow 180 44194 "Calameño" "Primera" 400 700 700 700 "Región de O'Higgins" 700
Set-DataRow 181 44194 "Calameño" "Segunda" 400 600 600 600 "Región de O'Higgins" 600
Set-DataRow 182 44194 "Tuna" "Extra" 400 800 800 800 "Región de O'Higgins" 800
Set-DataRow 183 44194 "Tuna" "Primera" 400 700 700 700 "Región de O'Higgins" 700
Set-DataRow 184 44194 "Tuna" "Segunda" 400 600 600 600 "Región de O'Higgins" 600

# Row 185: the "Calameño Primera" outlier record (previously row 179), date 44540.
Set-DataRow 185 44540 "Calameño" "Primera" 2700 1000 1100 1044 "Región de O'Higgins" 1044

# Rows 186-191: re-add the "Región Metropolitana" week (previously rows 180-185), date 44272.
Set-DataRow 186 44272 "Calameño" "Extra" 400 800 800 800 "Región Metropolitana" 800
Set-DataRow 187 44272 "Calameño" "Primera" 400 700 700 700 "Región Metropolitana" 700
Set-DataRow 188 44272 "Calameño" "Segunda" 400 500 500 500 "Región Metropolitana" 500
Set-DataRow 189 44272 "Tuna" "Extra" 400 800 800 800 "Región Metropolitana" 800
Set-DataRow 190 44272 "Tuna" "Primera" 400 700 700 700 "Región Metropolitana" 700
Set-DataRow 191 44272 "Tuna" "Segunda" 400 500 500 500 "Región Metropolitana" 500
